$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC!row 40
$wsALC.Cells.Item(40, 8).Value = 1553.375
$wsALC.Cells.Item(40, 9).Value = 883.3333
$wsALC.Cells.Item(40, 10).Value = 1955.4
$wsALC.Cells.Item(40, 11).Value = 883.3333
$wsALC.Cells.Item(40, 12).Value = 1955.4
$wsALC.Cells.Item(40, 13).Value = -708.3333
$wsALC.Cells.Item(40, 14).Value = -2305.4

# ALC!row 100
$wsALC.Cells.Item(100, 8).Value = 2045
$wsALC.Cells.Item(100, 9).Value = 1301.25
$wsALC.Cells.Item(100, 10).Value = 2640
$wsALC.Cells.Item(100, 11).Value = 1301.25
$wsALC.Cells.Item(100, 12).Value = 2640
$wsALC.Cells.Item(100, 13).Value = -760.25
$wsALC.Cells.Item(100, 14).Value = -3722

# ALC!row 129
$wsALC.Cells.Item(129, 8).Value = 748.60785
$wsALC.Cells.Item(129, 10).Value = 793.3043
$wsALC.Cells.Item(129, 12).Value = 2379.9129
$wsALC.Cells.Item(129, 14).Value = -12379.9129

# ALC!row 132
$wsALC.Cells.Item(132, 8).Value = 4273.9585
$wsALC.Cells.Item(132, 9).Value = 4438.8
$wsALC.Cells.Item(132, 10).Value = 3449.75
$wsALC.Cells.Item(132, 11).Value = 13316.4
$wsALC.Cells.Item(132, 12).Value = 10349.25
$wsALC.Cells.Item(132, 13).Value = -10786.4
$wsALC.Cells.Item(132, 14).Value = -15409.25

# ARM!row 2
$wsARM.Cells.Item(2, 8).Value = 1828.5714
$wsARM.Cells.Item(2, 9).Value = 1840
$wsARM.Cells.Item(2, 10).Value = 1800
$wsARM.Cells.Item(2, 11).Value = 1840
$wsARM.Cells.Item(2, 12).Value = 1800
$wsARM.Cells.Item(2, 13).Value = -1727
$wsARM.Cells.Item(2, 14).Value = -2026

# ARM!row 32
$wsARM.Cells.Item(32, 8).Value = 6440.1353
$wsARM.Cells.Item(32, 9).Value = 5195.053
$wsARM.Cells.Item(32, 10).Value = 10886.857
$wsARM.Cells.Item(32, 11).Value = 5195.053
$wsARM.Cells.Item(32, 12).Value = 10886.857
$wsARM.Cells.Item(32, 13).Value = -4908.053
$wsARM.Cells.Item(32, 14).Value = -11460.857

# ARM!row 116
$wsARM.Cells.Item(116, 8).Value = 1828.5714
$wsARM.Cells.Item(116, 9).Value = 1840
$wsARM.Cells.Item(116, 10).Value = 1800
$wsARM.Cells.Item(116, 11).Value = 1840
$wsARM.Cells.Item(116, 12).Value = 1800
$wsARM.Cells.Item(116, 13).Value = 454
$wsARM.Cells.Item(116, 14).Value = -6388

# ARM!row 132
$wsARM.Cells.Item(132, 8).Value = 26316.363
$wsARM.Cells.Item(132, 9).Value = 3579.5386
$wsARM.Cells.Item(132, 10).Value = 59158.445
$wsARM.Cells.Item(132, 11).Value = 10738.6158
$wsARM.Cells.Item(132, 12).Value = 177475.335
$wsARM.Cells.Item(132, 13).Value = -8208.6158
$wsARM.Cells.Item(132, 14).Value = -182535.335

# BSM!row 3
$wsBSM.Cells.Item(3, 8).Value = 1828.5714
$wsBSM.Cells.Item(3, 9).Value = 1840
$wsBSM.Cells.Item(3, 10).Value = 1800
$wsBSM.Cells.Item(3, 11).Value = 1840
$wsBSM.Cells.Item(3, 12).Value = 1800
$wsBSM.Cells.Item(3, 13).Value = -1726
$wsBSM.Cells.Item(3, 14).Value = -2028

# BSM!row 94
$wsBSM.Cells.Item(94, 8).Value = 1000.5
$wsBSM.Cells.Item(94, 9).Value = 850.6667
$wsBSM.Cells.Item(94, 10).Value = 1450
$wsBSM.Cells.Item(94, 11).Value = 850.6667
$wsBSM.Cells.Item(94, 12).Value = 1450
$wsBSM.Cells.Item(94, 13).Value = -399.6667
$wsBSM.Cells.Item(94, 14).Value = -2352

# CRP!row 31
$wsCRP.Cells.Item(31, 8).Value = 3544.383
$wsCRP.Cells.Item(31, 9).Value = 875.8946999999999
$wsCRP.Cells.Item(31, 10).Value = 5355.143
$wsCRP.Cells.Item(31, 11).Value = 875.8946999999999
$wsCRP.Cells.Item(31, 12).Value = 5355.143
$wsCRP.Cells.Item(31, 13).Value = -580.8946999999999
$wsCRP.Cells.Item(31, 14).Value = -5945.143

# CRP!row 34
$wsCRP.Cells.Item(34, 8).Value = 3544.383
$wsCRP.Cells.Item(34, 9).Value = 875.8946999999999
$wsCRP.Cells.Item(34, 10).Value = 5355.143
$wsCRP.Cells.Item(34, 11).Value = 875.8946999999999
$wsCRP.Cells.Item(34, 12).Value = 5355.143
$wsCRP.Cells.Item(34, 13).Value = -673.8946999999999
$wsCRP.Cells.Item(34, 14).Value = -5759.143

# CRP!row 99
$wsCRP.Cells.Item(99, 8).Value = 3770.3635
$wsCRP.Cells.Item(99, 9).Value = 2767.5293
$wsCRP.Cells.Item(99, 11).Value = 2767.5293
$wsCRP.Cells.Item(99, 13).Value = -1269.5293

# CRP!row 100
$wsCRP.Cells.Item(100, 8).Value = 122000
$wsCRP.Cells.Item(100, 10).Value = 122000
$wsCRP.Cells.Item(100, 12).Value = 122000
$wsCRP.Cells.Item(100, 14).Value = -124164

# CRP!row 126
$wsCRP.Cells.Item(126, 8).Value = 3770.3635
$wsCRP.Cells.Item(126, 9).Value = 2767.5293
$wsCRP.Cells.Item(126, 11).Value = 8302.5879
$wsCRP.Cells.Item(126, 13).Value = -5832.5879

# CRP!row 141
$wsCRP.Cells.Item(141, 8).Value = 33692.4
$wsCRP.Cells.Item(141, 10).Value = 33692.4
$wsCRP.Cells.Item(141, 12).Value = 33692.4
$wsCRP.Cells.Item(141, 14).Value = -44052.4

# CUL!row 5
$wsCUL.Cells.Item(5, 8).Value = 1101.6046
$wsCUL.Cells.Item(5, 9).Value = 758.871
$wsCUL.Cells.Item(5, 10).Value = 1987
$wsCUL.Cells.Item(5, 11).Value = 2276.613
$wsCUL.Cells.Item(5, 12).Value = 5961
$wsCUL.Cells.Item(5, 13).Value = -2164.613
$wsCUL.Cells.Item(5, 14).Value = -6185

# CUL!row 131
$wsCUL.Cells.Item(131, 8).Value = 763.5599999999999
$wsCUL.Cells.Item(131, 10).Value = 775.95874
$wsCUL.Cells.Item(131, 12).Value = 2327.87622
$wsCUL.Cells.Item(131, 14).Value = -12407.87622

# CUL!row 135
$wsCUL.Cells.Item(135, 8).Value = 1101.6046
$wsCUL.Cells.Item(135, 9).Value = 758.871
$wsCUL.Cells.Item(135, 10).Value = 1987
$wsCUL.Cells.Item(135, 11).Value = 6829.839
$wsCUL.Cells.Item(135, 12).Value = 17883
$wsCUL.Cells.Item(135, 13).Value = -4294.839
$wsCUL.Cells.Item(135, 14).Value = -22953

# GSM!row 15
$wsGSM.Cells.Item(15, 8).Value = 17500
$wsGSM.Cells.Item(15, 10).Value = 17500
$wsGSM.Cells.Item(15, 12).Value = 17500
$wsGSM.Cells.Item(15, 14).Value = -18076

# GSM!row 81
$wsGSM.Cells.Item(81, 8).Value = 17500
$wsGSM.Cells.Item(81, 10).Value = 17500
$wsGSM.Cells.Item(81, 12).Value = 17500
$wsGSM.Cells.Item(81, 14).Value = -19496

# GSM!row 84
$wsGSM.Cells.Item(84, 8).Value = 17500
$wsGSM.Cells.Item(84, 10).Value = 17500
$wsGSM.Cells.Item(84, 12).Value = 52500
$wsGSM.Cells.Item(84, 14).Value = -62484

# GSM!row 97
$wsGSM.Cells.Item(97, 8).Value = 1420.8572
$wsGSM.Cells.Item(97, 9).Value = 1555.9048
$wsGSM.Cells.Item(97, 10).Value = 1015.7143
$wsGSM.Cells.Item(97, 11).Value = 1555.9048
$wsGSM.Cells.Item(97, 12).Value = 1015.7143
$wsGSM.Cells.Item(97, 13).Value = -1059.9048
$wsGSM.Cells.Item(97, 14).Value = -2007.7143

# GSM!row 126
$wsGSM.Cells.Item(126, 8).Value = 2864.3147
$wsGSM.Cells.Item(126, 9).Value = 1879.963
$wsGSM.Cells.Item(126, 10).Value = 3848.6667
$wsGSM.Cells.Item(126, 11).Value = 5639.889
$wsGSM.Cells.Item(126, 12).Value = 11546.0001
$wsGSM.Cells.Item(126, 13).Value = -3169.889
$wsGSM.Cells.Item(126, 14).Value = -16486.0001

# LTW!row 47
$wsLTW.Cells.Item(47, 8).Value = 27495
$wsLTW.Cells.Item(47, 10).Value = 27495
$wsLTW.Cells.Item(47, 12).Value = 27495
$wsLTW.Cells.Item(47, 14).Value = -28475

# LTW!row 52
$wsLTW.Cells.Item(52, 8).Value = 27495
$wsLTW.Cells.Item(52, 10).Value = 27495
$wsLTW.Cells.Item(52, 12).Value = 27495
$wsLTW.Cells.Item(52, 14).Value = -27961

# LTW!row 100
$wsLTW.Cells.Item(100, 8).Value = 2329.2856
$wsLTW.Cells.Item(100, 9).Value = 1661
$wsLTW.Cells.Item(100, 10).Value = 2596.6
$wsLTW.Cells.Item(100, 11).Value = 1661
$wsLTW.Cells.Item(100, 12).Value = 2596.6
$wsLTW.Cells.Item(100, 13).Value = -1120
$wsLTW.Cells.Item(100, 14).Value = -3678.6

# LTW!row 122
$wsLTW.Cells.Item(122, 8).Value = 1035427.56
$wsLTW.Cells.Item(122, 9).Value = 1403887.4
$wsLTW.Cells.Item(122, 10).Value = 3740
$wsLTW.Cells.Item(122, 11).Value = 4211662.199999999
$wsLTW.Cells.Item(122, 12).Value = 11220
$wsLTW.Cells.Item(122, 13).Value = -4209212.199999999
$wsLTW.Cells.Item(122, 14).Value = -16120

# LTW!row 132
$wsLTW.Cells.Item(132, 8).Value = 672782.75
$wsLTW.Cells.Item(132, 9).Value = 862720.8
$wsLTW.Cells.Item(132, 10).Value = 7999.5
$wsLTW.Cells.Item(132, 11).Value = 2588162.4
$wsLTW.Cells.Item(132, 12).Value = 23998.5
$wsLTW.Cells.Item(132, 13).Value = -2585632.4
$wsLTW.Cells.Item(132, 14).Value = -29058.5

# WVR!row 75
$wsWVR.Cells.Item(75, 8).Value = 27500
$wsWVR.Cells.Item(75, 10).Value = 27500
$wsWVR.Cells.Item(75, 12).Value = 27500
$wsWVR.Cells.Item(75, 14).Value = -29372

# WVR!row 78
$wsWVR.Cells.Item(78, 8).Value = 27500
$wsWVR.Cells.Item(78, 10).Value = 27500
$wsWVR.Cells.Item(78, 12).Value = 82500
$wsWVR.Cells.Item(78, 14).Value = -91860

# WVR!row 132
$wsWVR.Cells.Item(132, 8).Value = 1251.4722
$wsWVR.Cells.Item(132, 9).Value = 679.5454999999999
$wsWVR.Cells.Item(132, 10).Value = 2150.2144
$wsWVR.Cells.Item(132, 11).Value = 2038.6365
$wsWVR.Cells.Item(132, 12).Value = 6450.6432
$wsWVR.Cells.Item(132, 13).Value = 491.3635000000002
$wsWVR.Cells.Item(132, 14).Value = -11510.6432
